{"js": "// Add a new block of meeting/interview notes at the end of the document body,\n// right after the existing \"10k + 20k + \" paragraph (CRUD: Create new rows\n// of notes for the \"Categories\" content).\nconst body = context.document.body;\n\nconst newParagraphs = [\n  \"Buying \",\n  \"Ibakuli 4k or 4.5k\",\n  \"Buying from farmers, then selling to the white man. Pays after 2-3 weeks, works with PALM (WFP), SAAWO, Canada. Buys from bikundi (farmers\\u2019 group \\u2013 beans, 2 hoes, forks) [umwiz\\u2019arahenda].\",\n  \"Other beans are now at 3k or 3.5k\",\n  \"Association is called Ajesedeki \\u2013 SCIATE (Works with this). $3 road construction. 580 (before), 450 (now) people.\",\n  \"Seasons: June & July | Jan & Feb\",\n  \"Ibakuli 5.5k\",\n  \"Arua Bus Park\",\n  \"\"\n];\n\nfor (const text of newParagraphs) {\n  body.insertParagraph(text, \"End\");\n}\n\nawait context.sync();\n", "ps1": "# Add a new block of meeting/interview notes at the end of the document body,\n# right after the existing \"10k + 20k + \" paragraph (CRUD: Create new rows\n# of notes for the \"Categories\" content).\n$d = $word.ActiveDocument\n\n$newParagraphs = @(\n    'Buying ',\n    'Ibakuli 4k or 4.5k',\n    'Buying from farmers, then selling to the white man. Pays after 2-3 weeks, works with PALM (WFP), SAAWO, Canada. Buys from bikundi (farmers\u2019 group \u2013 beans, 2 hoes, forks) [umwiz\u2019arahenda].',\n    'Other beans are now at 3k or 3.5k',\n    'Association is called Ajesedeki \u2013 SCIATE (Works with this). $3 road construction. 580 (before), 450 (now) people.',\n    'Seasons: June & July | Jan & Feb',\n    'Ibakuli 5.5k',\n    'Arua Bus Park',\n    ''\n)\n\nforeach ($t in $newParagraphs) {\n    $lastPara = $d.Paragraphs.Last\n    $r = $lastPara.Range\n    $r.InsertParagraphAfter()\n    $r.Collapse(0)\n    if ($t -ne '') {\n        $r.InsertAfter($t)\n    }\n}\n"}
